$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "50.897.19"
$ws.Range("E2").Value = "  -1.39%  "
Set-TextValue $ws.Range("D3") "2.931.74"
$ws.Range("E3").Value = "  +0.59%  "
Set-TextValue $ws.Range("D4") "0.995"
$ws.Range("E4").Value = "  -0.51%  "
Set-TextValue $ws.Range("D5") "370.62"
$ws.Range("E5").Value = "  +5.82%  "
Set-TextValue $ws.Range("D6") "104.43"
$ws.Range("E6").Value = "  -1.99%  "
Set-TextValue $ws.Range("D7") "0.541"
$ws.Range("E7").Value = "  -1.99%  "
Set-TextValue $ws.Range("D8") "1.03"
$ws.Range("E8").Value = "  +3.02%  "
Set-TextValue $ws.Range("D9") "0.589"
$ws.Range("E9").Value = "  -3.03%  "
Set-TextValue $ws.Range("D10") "36.82"
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("E11").Value = "  +1.26%  "
Set-TextValue $ws.Range("D12") "0.0837"
$ws.Range("E12").Value = "  -1.53%  "
Set-TextValue $ws.Range("D13") "18.40"
$ws.Range("E13").Value = "  -2.55%  "
Set-TextValue $ws.Range("D14") "3.378.35"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("E15").Value = "  -2.98%  "
Set-TextValue $ws.Range("D16") "2.911.26"
$ws.Range("E16").Value = "  -0.38%  "
Set-TextValue $ws.Range("D17") "0.941"
$ws.Range("E17").Value = "  -2.14%  "
Set-TextValue $ws.Range("D18") "50.902.52"
$ws.Range("E18").Value = "  -1.31%  "
Set-TextValue $ws.Range("D19") "3.27"
$ws.Range("E19").Value = "  -5.37%  "
Set-TextValue $ws.Range("D20") "7.24"
$ws.Range("E20").Value = "  -1.07%  "
Set-TextValue $ws.Range("D21") "12.95"
$ws.Range("E21").Value = "  -3.10%  "
Set-TextValue $ws.Range("D22") "0.0₃0946"
$ws.Range("E22").Value = "  -1.39%  "
Set-TextValue $ws.Range("D23") "68.32"
$ws.Range("E23").Value = "  -0.57%  "
Set-TextValue $ws.Range("D24") "260.47"
$ws.Range("E24").Value = "  -0.16%  "
Set-TextValue $ws.Range("D25") "2.69"
$ws.Range("E25").Value = "  +0.22%  "
Set-TextValue $ws.Range("D26") "4.33"
$ws.Range("E26").Value = "  +3.67%  "
Set-TextValue $ws.Range("D27") "0.173"
$ws.Range("E27").Value = "  +1.34%  "
Set-TextValue $ws.Range("D28") "0.999"
$ws.Range("E28").Value = "  -0.05%  "
Set-TextValue $ws.Range("D29") "25.83"
$ws.Range("E29").Value = "  -2.18%  "
Set-TextValue $ws.Range("D30") "7.01"
$ws.Range("E30").Value = "  -8.15%  "
Set-TextValue $ws.Range("D31") "0.104"
$ws.Range("E31").Value = "  +0.53%  "
Set-TextValue $ws.Range("D32") "6.25"
$ws.Range("E32").Value = "  +1.81%  "
Set-TextValue $ws.Range("D33") "9.92"
$ws.Range("E33").Value = "  -2.52%  "
$ws.Range("E34").Value = "  -1.42%  "
Set-TextValue $ws.Range("D35") "34.91"
$ws.Range("E35").Value = "  -1.60%  "
Set-TextValue $ws.Range("D36") "50.90"
$ws.Range("E36").Value = "  +0.22%  "
Set-TextValue $ws.Range("D37") "0.998"
$ws.Range("E37").Value = "  -0.16%  "
Set-TextValue $ws.Range("D38") "0.0423"
$ws.Range("E38").Value = "  +0.05%  "
Set-TextValue $ws.Range("D39") "3.05"
$ws.Range("E39").Value = "  -1.70%  "
Set-TextValue $ws.Range("D40") "2.67"
$ws.Range("E40").Value = "  +1.46%  "
Set-TextValue $ws.Range("D41") "17.13"
$ws.Range("E41").Value = "  -2.25%  "
Set-TextValue $ws.Range("D42") "1.85"
$ws.Range("E42").Value = "  -4.26%  "
$ws.Range("E43").Value = "  -1.68%  "
Set-TextValue $ws.Range("D44") "22.15"
$ws.Range("E44").Value = "  -0.81%  "
Set-TextValue $ws.Range("D45") "119.33"
$ws.Range("E45").Value = "  +0.04%  "
Set-TextValue $ws.Range("D46") "2.08"
$ws.Range("E46").Value = "  -2.58%  "
Set-TextValue $ws.Range("D47") "2.027.08"
$ws.Range("E47").Value = "  -2.90%  "
$ws.Range("E48").Value = "  -0.16%  "
Set-TextValue $ws.Range("D49") "3.19"
$ws.Range("E49").Value = "  -3.41%  "
Set-TextValue $ws.Range("D50") "0.247"
$ws.Range("E50").Value = "  +4.60%  "
Set-TextValue $ws.Range("D51") "3.202.58"
$ws.Range("E51").Value = "  -0.30%  "
